# Update NATMI LR-pair TPM statistics (Efnb2-Ephb3) with recomputed values
# from the refreshed TPM pipeline run ("update scripts wuth new tpm").
# Columns E:T across data rows 2-17 are refreshed to the newly computed
# ligand/receptor expression, detection, and specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.76217133333333
$ws.Range("H2").Value = 137.286514
$ws.Range("I2").Value = 0.6763939203605134
$ws.Range("J2").Value = 0.6763939203605135
$ws.Range("M2").Value = 0.103005
$ws.Range("N2").Value = 0.309015
$ws.Range("O2").Value = 0.004680825815734043
$ws.Range("P2").Value = 0.004680825815734042
$ws.Range("Q2").Value = 4.71373245819
$ws.Range("R2").Value = 42.42359212371
$ws.Range("S2").Value = 0.003166082124029048
$ws.Range("T2").Value = 0.003166082124029048
$ws.Range("G3").Value = 45.76217133333333
$ws.Range("H3").Value = 137.286514
$ws.Range("I3").Value = 0.6763939203605134
$ws.Range("J3").Value = 0.6763939203605135
$ws.Range("N3").Value = 54.36902
$ws.Range("O3").Value = 0.8235584434158876
$ws.Range("P3").Value = 0.8235584434158876
$ws.Range("Q3").Value = 829.3481361551421
$ws.Range("R3").Value = 7464.13322539628
$ws.Range("S3").Value = 0.5570499241880742
$ws.Range("T3").Value = 0.5570499241880743
$ws.Range("G4").Value = 45.76217133333333
$ws.Range("H4").Value = 137.286514
$ws.Range("I4").Value = 0.6763939203605134
$ws.Range("J4").Value = 0.6763939203605135
$ws.Range("M4").Value = 3.761887333333334
$ws.Range("N4").Value = 11.285662
$ws.Range("O4").Value = 0.1709503358647596
$ws.Range("P4").Value = 0.1709503358647596
$ws.Range("Q4").Value = 172.1521326846965
$ws.Range("R4").Value = 1549.369194162268
$ws.Range("S4").Value = 0.1156297678625112
$ws.Range("T4").Value = 0.1156297678625112
$ws.Range("G5").Value = 45.76217133333333
$ws.Range("H5").Value = 137.286514
$ws.Range("I5").Value = 0.6763939203605134
$ws.Range("J5").Value = 0.6763939203605135
$ws.Range("M5").Value = 0.01783333333333333
$ws.Range("N5").Value = 0.0535
$ws.Range("O5").Value = 0.0008103949036188254
$ws.Range("P5").Value = 0.0008103949036188253
$ws.Range("Q5").Value = 0.8160920554444444
$ws.Range("R5").Value = 7.344828499
$ws.Range("S5").Value = 0.0005481461858989177
$ws.Range("T5").Value = 0.0005481461858989177
$ws.Range("I6").Value = 0.1388778842960613
$ws.Range("J6").Value = 0.1388778842960613
$ws.Range("M6").Value = 0.103005
$ws.Range("N6").Value = 0.309015
$ws.Range("O6").Value = 0.004680825815734043
$ws.Range("P6").Value = 0.004680825815734042
$ws.Range("Q6").Value = 0.9678283190099999
$ws.Range("R6").Value = 8.710454871089999
$ws.Range("S6").Value = 0.0006500631860475291
$ws.Range("T6").Value = 0.0006500631860475291
$ws.Range("I7").Value = 0.1388778842960613
$ws.Range("J7").Value = 0.1388778842960613
$ws.Range("N7").Value = 54.36902
$ws.Range("O7").Value = 0.8235584434158876
$ws.Range("P7").Value = 0.8235584434158876
$ws.Range("S7").Value = 0.114374054215756
$ws.Range("T7").Value = 0.114374054215756
$ws.Range("I8").Value = 0.1388778842960613
$ws.Range("J8").Value = 0.1388778842960613
$ws.Range("M8").Value = 3.761887333333334
$ws.Range("N8").Value = 11.285662
$ws.Range("O8").Value = 0.1709503358647596
$ws.Range("P8").Value = 0.1709503358647596
$ws.Range("Q8").Value = 35.34645011528578
$ws.Range("R8").Value = 318.118051037572
$ws.Range("S8").Value = 0.0237412209645989
$ws.Range("T8").Value = 0.02374122096459891
$ws.Range("I9").Value = 0.1388778842960613
$ws.Range("J9").Value = 0.1388778842960613
$ws.Range("M9").Value = 0.01783333333333333
$ws.Range("N9").Value = 0.0535
$ws.Range("O9").Value = 0.0008103949036188254
$ws.Range("P9").Value = 0.0008103949036188253
$ws.Range("Q9").Value = 0.1675608467777777
$ws.Range("R9").Value = 1.508047621
$ws.Range("S9").Value = 0.000112545929658893
$ws.Range("T9").Value = 0.000112545929658893
$ws.Range("G10").Value = 12.29750866666667
$ws.Range("H10").Value = 36.892526
$ws.Range("I10").Value = 0.1817649787009828
$ws.Range("J10").Value = 0.1817649787009828
$ws.Range("M10").Value = 0.103005
$ws.Range("N10").Value = 0.309015
$ws.Range("O10").Value = 0.004680825815734043
$ws.Range("P10").Value = 0.004680825815734042
$ws.Range("Q10").Value = 1.26670488021
$ws.Range("R10").Value = 11.40034392189
$ws.Range("S10").Value = 0.0008508102046999086
$ws.Range("T10").Value = 0.0008508102046999085
$ws.Range("G11").Value = 12.29750866666667
$ws.Range("H11").Value = 36.892526
$ws.Range("I11").Value = 0.1817649787009828
$ws.Range("J11").Value = 0.1817649787009828
$ws.Range("N11").Value = 54.36902
$ws.Range("O11").Value = 0.8235584434158876
$ws.Range("P11").Value = 0.8235584434158876
$ws.Range("Q11").Value = 222.8678315493911
$ws.Range("R11").Value = 2005.81048394452
$ws.Range("S11").Value = 0.1496940829265033
$ws.Range("T11").Value = 0.1496940829265033
$ws.Range("G12").Value = 12.29750866666667
$ws.Range("H12").Value = 36.892526
$ws.Range("I12").Value = 0.1817649787009828
$ws.Range("J12").Value = 0.1817649787009828
$ws.Range("M12").Value = 3.761887333333334
$ws.Range("N12").Value = 11.285662
$ws.Range("O12").Value = 0.1709503358647596
$ws.Range("P12").Value = 0.1709503358647596
$ws.Range("Q12").Value = 46.26184208469023
$ws.Range("R12").Value = 416.3565787622121
$ws.Range("S12").Value = 0.03107278415738389
$ws.Range("T12").Value = 0.03107278415738388
$ws.Range("G13").Value = 12.29750866666667
$ws.Range("H13").Value = 36.892526
$ws.Range("I13").Value = 0.1817649787009828
$ws.Range("J13").Value = 0.1817649787009828
$ws.Range("M13").Value = 0.01783333333333333
$ws.Range("N13").Value = 0.0535
$ws.Range("O13").Value = 0.0008103949036188254
$ws.Range("P13").Value = 0.0008103949036188253
$ws.Range("Q13").Value = 0.2193055712222222
$ws.Range("R13").Value = 1.973750141
$ws.Range("S13").Value = 0.0001473014123956608
$ws.Range("T13").Value = 0.0001473014123956607
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2004796666666666
$ws.Range("H14").Value = 0.6014389999999999
$ws.Range("I14").Value = 0.002963216642442438
$ws.Range("J14").Value = 0.002963216642442439
$ws.Range("M14").Value = 0.103005
$ws.Range("N14").Value = 0.309015
$ws.Range("O14").Value = 0.004680825815734043
$ws.Range("P14").Value = 0.004680825815734042
$ws.Range("Q14").Value = 0.020650408065
$ws.Range("R14").Value = 0.185853672585
$ws.Range("S14").Value = 0.00001387030095755732
$ws.Range("T14").Value = 0.00001387030095755732
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2004796666666666
$ws.Range("H15").Value = 0.6014389999999999
$ws.Range("I15").Value = 0.002963216642442438
$ws.Range("J15").Value = 0.002963216642442439
$ws.Range("N15").Value = 54.36902
$ws.Range("O15").Value = 0.8235584434158876
$ws.Range("P15").Value = 0.8235584434158876
$ws.Range("Q15").Value = 3.63329433553111
$ws.Range("R15").Value = 32.69964901978
$ws.Range("S15").Value = 0.002440382085553947
$ws.Range("T15").Value = 0.002440382085553947
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2004796666666666
$ws.Range("H16").Value = 0.6014389999999999
$ws.Range("I16").Value = 0.002963216642442438
$ws.Range("J16").Value = 0.002963216642442439
$ws.Range("M16").Value = 3.761887333333334
$ws.Range("N16").Value = 11.285662
$ws.Range("O16").Value = 0.1709503358647596
$ws.Range("P16").Value = 0.1709503358647596
$ws.Range("Q16").Value = 0.7541819186242221
$ws.Range("R16").Value = 6.787637267618
$ws.Range("S16").Value = 0.0005065628802655802
$ws.Range("T16").Value = 0.0005065628802655802
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2004796666666666
$ws.Range("H17").Value = 0.6014389999999999
$ws.Range("I17").Value = 0.002963216642442438
$ws.Range("J17").Value = 0.002963216642442439
$ws.Range("M17").Value = 0.01783333333333333
$ws.Range("N17").Value = 0.0535
$ws.Range("O17").Value = 0.0008103949036188254
$ws.Range("P17").Value = 0.0008103949036188253
$ws.Range("Q17").Value = 0.003575220722222222
$ws.Range("R17").Value = 0.0321769865
$ws.Range("S17").Value = 0.000002401375665353839
$ws.Range("T17").Value = 0.000002401375665353839

Write-Host "Applied updated TPM values to Efnb2-Ephb3 LR-pair sheet."
